$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for "Mã cán bộ" (shift existing A:C to B:D)
$ws.Columns("A:A").Insert()

# Header row (row-major order so new shared strings are appended in
# reading order: A1, B1, C1, D1, A2, B2, ...)
$ws.Range("A1").Value = "Mã cán bộ"
$ws.Range("B1").Value = "Họ tên"
$ws.Range("C1").Value = "Hình ảnh"
$ws.Range("D1").Value = "Tên bộ phận"

# Row 2
$ws.Range("A2").Value = "cb14"
$ws.Range("B2").Value = "Trần Quốc K"
$ws.Range("C2").Value = "D:\meme.png"
$ws.Range("D2").Value = "Lĩnh vực lao động, thương binh và xã hội"

# Row 3
$ws.Range("A3").Value = "cb15"
$ws.Range("B3").Value = "Trần Trung P"
$ws.Range("C3").Value = "D:\meme.png"
$ws.Range("D3").Value = "Lĩnh vực lao động, Thương Binh Và Xã Hội"

# Update selection / view state
$ws.Range("F7").Select()
